$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-19T08:22:07+00:00 -> 2025-12-19T09:47:21+00:00
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# Base Definition: add version suffix |4.0.1
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s) cell (K6): add version pins to both referenced profiles
$elem.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-practitioner|2.2.0-ballot|https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-practitioner-role|2.2.0-ballot)`n"

# Column K width grew to fit the longer text (raw OOXML width 139.5390625 -> 158.5546875).
# ColumnWidth is expressed in character units; the stored raw width equals ColumnWidth + 5/6.
$elem.Columns.Item(11).ColumnWidth = 158.5546875 - (5/6)
